$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values - use quote-prefix to force text storage
# so values like "4.60" or "1.00" are not coerced into numbers and lose formatting.
$ws.Range("D2").Value = "'68.813.40"
$ws.Range("D3").Value = "'2.442.45"
$ws.Range("D5").Value = "'560.72"
$ws.Range("D6").Value = "'163.67"
$ws.Range("D8").Value = "'0.507"
$ws.Range("D10").Value = "'0.161"
$ws.Range("D12").Value = "'4.60"
$ws.Range("D14").Value = "'68.699.19"
$ws.Range("D15").Value = "'2.889.99"
$ws.Range("D16").Value = "'23.43"
$ws.Range("D17").Value = "'2.442.89"
$ws.Range("D18").Value = "'10.63"
$ws.Range("D19").Value = "'339.07"
$ws.Range("D20").Value = "'7.02"
$ws.Range("D22").Value = "'1.94"
$ws.Range("D24").Value = "'65.41"
$ws.Range("D26").Value = "'2.568.43"
$ws.Range("D27").Value = "'8.40"
$ws.Range("D28").Value = "'1.01"
$ws.Range("D29").Value = "'0.0₃0826"
$ws.Range("D31").Value = "'1.21"
$ws.Range("D32").Value = "'1.00"
$ws.Range("D33").Value = "'433.35"
$ws.Range("D35").Value = "'159.33"
$ws.Range("D38").Value = "'18.01"
$ws.Range("D44").Value = "'2.09"
$ws.Range("D46").Value = "'129.94"
$ws.Range("D47").Value = "'0.0719"
$ws.Range("D48").Value = "'0.485"

# Update Volume(1h) (column E) values
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.81%  "
$ws.Range("E9").Value = "  +8.99%  "
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("E12").Value = "  -4.92%  "
$ws.Range("E13").Value = "  +5.04%  "
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("E22").Value = "  +2.29%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -2.06%  "
$ws.Range("E25").Value = "  +2.90%  "
$ws.Range("E26").Value = "  -0.94%  "
$ws.Range("E27").Value = "  +2.51%  "
$ws.Range("E28").Value = "  +1.29%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("E31").Value = "  +5.17%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("E34").Value = "  -1.83%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("E40").Value = "  +1.10%  "
$ws.Range("E41").Value = "  +2.28%  "
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("E44").Value = "  +1.38%  "
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("E51").Value = "  +2.51%  "
